# Update "想去人数" (interested-count) figures to the latest scraped values.
$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1535
$ws1.Range("F6").Value  = 2309
$ws1.Range("F11").Value = 350
$ws1.Range("F12").Value = 1094
$ws1.Range("F17").Value = 4709
$ws1.Range("F19").Value = 1333
$ws1.Range("F20").Value = 3508
$ws1.Range("F24").Value = 3718
$ws1.Range("F25").Value = 5119
$ws1.Range("F27").Value = 980
$ws1.Range("F28").Value = 562
$ws1.Range("F39").Value = 137
$ws1.Range("F40").Value = 1383
$ws1.Range("F42").Value = 858
$ws1.Range("F45").Value = 341
$ws1.Range("F47").Value = 164
$ws1.Range("F49").Value = 3739

# --- 演出 (Performances) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 1016

# --- 本地生活 (Local life) sheet ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2277

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 2277
$ws4.Range("F3").Value  = 1535
$ws4.Range("F8").Value  = 2309
$ws4.Range("F12").Value = 1016
$ws4.Range("F15").Value = 350
$ws4.Range("F16").Value = 1094
$ws4.Range("F21").Value = 4709
$ws4.Range("F22").Value = 1333
$ws4.Range("F24").Value = 3508
$ws4.Range("F25").Value = 3718
$ws4.Range("F26").Value = 5119
$ws4.Range("F28").Value = 980
$ws4.Range("F29").Value = 562
$ws4.Range("F38").Value = 137
$ws4.Range("F39").Value = 1383
$ws4.Range("F45").Value = 341
$ws4.Range("F47").Value = 164
$ws4.Range("F49").Value = 3739
